# Update the Extent report workbook to reflect the latest test run
# (new timings / scenario-feature names / tag / duration numbers),
# mirroring what the ExtentReports library writes on a fresh run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# DB Data sheet (very-hidden) - run summary block
# ---------------------------------------------------------------
$db = $wb.Worksheets.Item("DB Data")
$db.Range("B3").Value = "Jan 01, 2024 5:11:44 PM"   # Date Time
$db.Range("B4").Value = "Jan 01, 2024 5:09:34 PM"   # Start Date Time
$db.Range("B5").Value = "Jan 01, 2024 5:11:38 PM"   # End Date Time
$db.Range("B6").Value = "2 m 4.207 s"                # Duration
$db.Range("H2").Value = 18                            # Scenarios Total (feeds H5 SUM)

# ---------------------------------------------------------------
# Scenarios sheet - single scenario row
# ---------------------------------------------------------------
$scenarios = $wb.Worksheets.Item("Scenarios")
$scenarios.Range("B22").Value = "User should verify address details in checkout page"
$scenarios.Range("D22").Value = "2 m 3.209 s"
$scenarios.Range("E22").Value = "Register Feature"
$scenarios.Range("G22").Value = 18
$scenarios.Range("H22").Value = 18

# ---------------------------------------------------------------
# Features sheet - single feature row
# ---------------------------------------------------------------
$features = $wb.Worksheets.Item("Features")
$features.Range("B22").Value = "Register Feature"
$features.Range("D22").Value = "2 m 3.215 s"
$features.Range("J22").Value = 18
$features.Range("K22").Value = 18

# ---------------------------------------------------------------
# Tags sheet - tag name + feature/scenario references
# ---------------------------------------------------------------
$tags = $wb.Worksheets.Item("Tags")
$tags.Range("B24").Value = "@register"
$tags.Range("C29").Value = "Register Feature"
$tags.Range("C30").Value = "Register Feature"
$tags.Range("C31").Value = "Register Feature"
$tags.Range("H29").Value = "User should verify address details in checkout page"
$tags.Range("H30").Value = "User should verify address details in checkout page"
$tags.Range("H31").Value = "User should verify address details in checkout page"
